# Replace the four "Campaign Dates that use Perseus: ..." paragraphs with a
# single, plain-formatted run reading the translated Bootes dates, matching
# the author's commit that collapses the old multi-run date string down to
# one run with no explicit run properties.

$d = $word.ActiveDocument

$newText = "Campaign Dates that use Bootes: May 14-23, June 13-22, July 12-21"

# A minimal single-part WordprocessingML package whose body is just the one
# paragraph/run we want. InsertXML() replaces the target Range's contents
# with this, which (unlike setting .Text) lets us land a run with no <w:rPr>
# at all -- exactly what the target markup has.
$xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Locate every paragraph still carrying the old campaign-dates sentence.
$targets = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Campaign Dates that use Perseus*") {
        $targets += $i
    }
}

# Walk back-to-front so replacing one paragraph never disturbs the index of
# another we still need to visit.
for ($j = $targets.Count - 1; $j -ge 0; $j--) {
    $p = $d.Paragraphs.Item($targets[$j])
    $full = $p.Range
    # Exclude the trailing paragraph mark from the range so InsertXML only
    # rewrites the runs, leaving the paragraph (and its pPr) itself intact.
    $body = $d.Range($full.Start, $full.End - 1)
    $body.InsertXML($xmlFrag)
}
